$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Nikola"
$ws.Range("B8").Value = "Jokic"
$ws.Range("C8").Value = 11
$ws.Range("A9").Value = "Ben"
$ws.Range("B9").Value = "Wallace"
$ws.Range("C9").Value = 10
$ws.Range("A10").Value = "Marcus"
$ws.Range("B10").Value = "Camby"
$ws.Range("C10").Value = 9
$ws.Range("A11").Value = "Lafayette"
$ws.Range("B11").Value = "Lever"
$ws.Range("C11").Value = 9
$ws.Range("A12").Value = "Anthony"
$ws.Range("B12").Value = "Davis"
$ws.Range("C12").Value = 9
$ws.Range("A29").Value = "Larry"
$ws.Range("B29").Value = "Steele"
$ws.Range("C29").Value = 4
$ws.Range("A34").Value = "Luka"
$ws.Range("B34").Value = "Doncic"
$ws.Range("C34").Value = 4
$ws.Range("A36").Value = "Mark"
$ws.Range("B36").Value = "Eaton"
$ws.Range("C36").Value = 4
$ws.Range("A37").Value = "Gerald"
$ws.Range("B37").Value = "Wallace"
$ws.Range("C37").Value = 4
$ws.Range("A38").Value = "Giannis"
$ws.Range("B38").Value = "Antetokounmpo"
$ws.Range("C38").Value = 4
$ws.Range("A39").Value = "Bob"
$ws.Range("B39").Value = "McAdoo"
$ws.Range("C39").Value = 4
$ws.Range("A40").Value = "Clyde"
$ws.Range("B40").Value = "Drexler"
$ws.Range("C40").Value = 4
$ws.Range("A42").Value = "George"
$ws.Range("B42").Value = "McGinnis"
$ws.Range("C42").Value = 3
$ws.Range("A89").Value = "Charles"
$ws.Range("B89").Value = "Jones"
$ws.Range("C89").Value = 2
$ws.Range("A90").Value = "Damian"
$ws.Range("B90").Value = "Lillard"
$ws.Range("C90").Value = 2
$ws.Range("A101").Value = "Paul"
$ws.Range("B101").Value = "Pressey"
$ws.Range("C101").Value = 1
$ws.Range("A102").Value = "Quinn"
$ws.Range("B102").Value = "Buckner"
$ws.Range("C102").Value = 1
$ws.Range("A107").Value = "Brook"
$ws.Range("B107").Value = "Lopez"
$ws.Range("C107").Value = 1
$ws.Range("A112").Value = "Nicolas"
$ws.Range("B112").Value = "Batum"
$ws.Range("C112").Value = 1
$ws.Range("A113").Value = "Nick"
$ws.Range("B113").Value = "Anderson"
$ws.Range("C113").Value = 1
$ws.Range("A114").Value = "Cade"
$ws.Range("B114").Value = "Cunningham"
$ws.Range("C114").Value = 1
$ws.Range("A115").Value = "Charles"
$ws.Range("B115").Value = "Oakley"
$ws.Range("C115").Value = 1
$ws.Range("A116").Value = "Charles"
$ws.Range("B116").Value = "Smith"
$ws.Range("C116").Value = 1
$ws.Range("A117").Value = "Moussa"
$ws.Range("B117").Value = "Diabate"
$ws.Range("C117").Value = 1
